# Update cryptocurrency price column (D) with refreshed quotes.
# Values are stored as text (inline strings) in the sheet, not numbers, so
# assigning a bare numeric-looking string would make Excel auto-convert the
# cell to a floating point number (losing exact text such as trailing
# zeros). Prefixing with a literal single-quote forces Excel to keep the
# entry as text; ClearFormats() afterwards strips the "quote prefix" cell
# style that operation leaves behind, restoring the cell to its original,
# unstyled state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = "264.72"
    3  = "22.72"
    4  = "6.230"
    5  = "0.06127"
    6  = "3.550"
    8  = "1.374"
    9  = "0.8139"
    10 = "0.1596"
    11 = "0.08236"
    12 = "0.03378"
    14 = "0.09263"
    15 = "3.914"
    16 = "0.001702"
    17 = "0.04848"
    18 = "0.0006275"
    19 = "0.006249"
    20 = "0.001106"
    21 = "0.003203"
    22 = "0.0001507"
    23 = "3.693"
    24 = "2.258"
    25 = "0.3391"
    26 = "0.1272"
    40 = "0.04630"
    41 = "0.007292"
    42 = "0.1126"
    43 = "0.003403"
    45 = "0.00006184"
    47 = "0.7531"
    48 = "0.1952"
    49 = "0.00002109"
    50 = "0.01245"
}

foreach ($row in $updates.Keys) {
    $cell = $ws.Range("D$row")
    $cell.Value = "'" + $updates[$row]
    $cell.ClearFormats()
}
